# Atualização de bases das ligas, do dia: 14-04-2024 às 15:12
#
# This script updates match rows 119/120 (odds for two already-completed
# matches get corrected/swapped back into the right rows) and refreshes the
# upcoming-fixture block at rows 177-182: several matches get their closing
# odds refreshed, one old (now irrelevant) fixture is dropped, new fixtures
# are introduced, and the table grows by one row (182).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 119: id 117 - Nacional De Football vs Torque (was showing the
# Defensor Sporting vs Danubio data; corrected to its own data)
# ---------------------------------------------------------------------
$ws.Range("B119").Value = 7013409
$ws.Range("F119").Value = "Nacional De Football"
$ws.Range("G119").Value = "Torque"
$ws.Range("H119").Value = 1
$ws.Range("I119").Value = 1
$ws.Range("J119").Value = "D"
$ws.Range("K119").Value = 1.666
$ws.Range("L119").Value = 3.9
$ws.Range("M119").Value = 4.5
$ws.Range("N119").Value = 1.615
$ws.Range("O119").Value = 4
$ws.Range("P119").Value = 4.75
$ws.Range("Q119").Value = -0.75
$ws.Range("R119").Value = 1.8
$ws.Range("S119").Value = 2.05
$ws.Range("T119").Value = 2.75
$ws.Range("U119").Value = 1.95
$ws.Range("V119").Value = 1.9
$ws.Range("W119").Value = -1
$ws.Range("X119").Value = 3
$ws.Range("Y119").Value = -1
$ws.Range("Z119").Value = -1
$ws.Range("AA119").Value = 1.05
$ws.Range("AB119").Value = -1
$ws.Range("AC119").Value = 0.8999999999999999

# ---------------------------------------------------------------------
# Row 120: id 118 - Defensor Sporting vs Danubio (the data that used to
# be in row 119)
# ---------------------------------------------------------------------
$ws.Range("B120").Value = 7013702
$ws.Range("F120").Value = "Defensor Sporting"
$ws.Range("G120").Value = "Danubio"
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 2
$ws.Range("J120").Value = "A"
$ws.Range("K120").Value = 1.8
$ws.Range("L120").Value = 3.6
$ws.Range("M120").Value = 4.2
$ws.Range("N120").Value = 1.8
$ws.Range("O120").Value = 3.6
$ws.Range("P120").Value = 4.2
$ws.Range("Q120").Value = -0.75
$ws.Range("R120").Value = 2.05
$ws.Range("S120").Value = 1.8
$ws.Range("T120").Value = 2.25
$ws.Range("U120").Value = 1.85
$ws.Range("V120").Value = 2
$ws.Range("W120").Value = -1
$ws.Range("X120").Value = -1
$ws.Range("Y120").Value = 3.2
$ws.Range("Z120").Value = -1
$ws.Range("AA120").Value = 0.8
$ws.Range("AB120").Value = -0.5
$ws.Range("AC120").Value = 0.5

# ---------------------------------------------------------------------
# Row 177: id 175 dropped; row 176 (Cerro vs Atletico Fenix Montevideo,
# id 8051185) shifts up with refreshed closing odds
# ---------------------------------------------------------------------
$ws.Range("B177").Value = 8051185
$ws.Range("E177").Value = 45395.72916666666
$ws.Range("F177").Value = "Cerro"
$ws.Range("G177").Value = "Atletico Fenix Montevideo"
$ws.Range("L177").Value = 3
$ws.Range("M177").Value = 3.2
$ws.Range("N177").Value = 2.6
$ws.Range("O177").Value = 3
$ws.Range("P177").Value = 2.9
$ws.Range("R177").Value = 1.8
$ws.Range("S177").Value = 2.05
$ws.Range("T177").Value = 2
$ws.Range("U177").Value = 1.85
$ws.Range("V177").Value = 2

# ---------------------------------------------------------------------
# Row 178: Deportivo Maldonado vs Cerro Largo, id 8051186, shifts up
# with refreshed closing odds
# ---------------------------------------------------------------------
$ws.Range("B178").Value = 8051186
$ws.Range("E178").Value = 45395.83333333334
$ws.Range("F178").Value = "Deportivo Maldonado"
$ws.Range("G178").Value = "Cerro Largo"
$ws.Range("K178").Value = 2.3
$ws.Range("N178").Value = 2.375
$ws.Range("P178").Value = 3.2
$ws.Range("Q178").Value = -0.25
$ws.Range("R178").Value = 2.05
$ws.Range("S178").Value = 1.8
$ws.Range("T178").Value = 2.25
$ws.Range("U178").Value = 2
$ws.Range("V178").Value = 1.85

# ---------------------------------------------------------------------
# Row 179: brand-new fixture - Rampla Juniors vs Racing Club de
# Montevideo, id 8051004 (old id-8051186 data that used to sit here has
# already moved to row 178)
# ---------------------------------------------------------------------
$ws.Range("B179").Value = 8051004
$ws.Range("E179").Value = 45396.41666666666
$ws.Range("F179").Value = "Rampla Juniors"
$ws.Range("G179").Value = "Racing Club de Montevideo"
$ws.Range("K179").Value = 3.2
$ws.Range("L179").Value = 3.3
$ws.Range("M179").Value = 2.2
$ws.Range("N179").Value = 3.4
$ws.Range("O179").Value = 3.5
$ws.Range("P179").Value = 2.05
$ws.Range("Q179").Value = 0.25
$ws.Range("T179").Value = 2.5
$ws.Range("U179").Value = 2.025
$ws.Range("V179").Value = 1.825

# Row 180 (id 178, Defensor Sporting vs CA River Plate) is unchanged.

# ---------------------------------------------------------------------
# Row 181: brand-new fixture - Penarol vs Danubio, id 8050911
# ---------------------------------------------------------------------
$ws.Range("B181").Value = 8050911
$ws.Range("E181").Value = 45396.75
$ws.Range("F181").Value = "Penarol"
$ws.Range("G181").Value = "Danubio"
$ws.Range("K181").Value = 1.666
$ws.Range("L181").Value = 3.5
$ws.Range("M181").Value = 5.5
$ws.Range("N181").Value = 1.65
$ws.Range("O181").Value = 3.5
$ws.Range("P181").Value = 5.75
$ws.Range("Q181").Value = -0.75
$ws.Range("R181").Value = 1.825
$ws.Range("S181").Value = 2.025
$ws.Range("U181").Value = 1.925
$ws.Range("V181").Value = 1.925

# ---------------------------------------------------------------------
# Row 182 (new row): the fixture that used to be row 181 - Montevideo
# Wanderers vs Liverpool Montevideo, id 8050912 - moves down one row.
# Copy formatting (bold/bordered id style, date format) from existing
# rows first, then fill in the values.
# ---------------------------------------------------------------------
$ws.Range("A119").Copy($ws.Range("A182"))
$ws.Range("E119").Copy($ws.Range("E182"))

$ws.Range("A182").Value = 180
$ws.Range("B182").Value = 8050912
$ws.Range("C182").Value = "Uruguay Primera División"
$ws.Range("D182").Value = "Uruguay Apertura"
$ws.Range("E182").Value = 45397.75
$ws.Range("F182").Value = "Montevideo Wanderers"
$ws.Range("G182").Value = "Liverpool Montevideo"
$ws.Range("K182").Value = 3.2
$ws.Range("L182").Value = 3.3
$ws.Range("M182").Value = 2.2
$ws.Range("N182").Value = 3.4
$ws.Range("O182").Value = 3.3
$ws.Range("P182").Value = 2.1
$ws.Range("Q182").Value = 0.25
$ws.Range("R182").Value = 2
$ws.Range("S182").Value = 1.85
$ws.Range("T182").Value = 2.25
$ws.Range("U182").Value = 1.9
$ws.Range("V182").Value = 1.95
$ws.Range("W182").Value = 0
$ws.Range("X182").Value = 0
$ws.Range("Y182").Value = 0
$ws.Range("Z182").Value = 0
$ws.Range("AA182").Value = 0
